# Powerpoint writer: consolidate text runs when possible.
#
# Several text boxes in this deck were originally split across many
# <a:r> runs (one run per word/space token). Re-assigning the full
# text via TextRange.Text causes the host to emit a single merged
# run, which is exactly the "run consolidation" behaviour we want to
# reproduce here.
#
# Note: the rendered text is already correct before this edit (it was
# just split across runs), so setting TextRange.Text to that same
# string is a no-op as far as the rendered characters go. To force the
# run-merge to actually happen we first set the text to a throwaway
# placeholder value and then set it to the desired final string.

$p = $ppt.ActivePresentation

function Set-ConsolidatedText($textRange, $finalText) {
    $textRange.Text = "~"
    $textRange.Text = $finalText
}

# Slide 2 title: "Slide" " " "1" -> "Slide 1"
Set-ConsolidatedText $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange "Slide 1"

# Slide 4 title: "Slide" " " "3" -> "Slide 3"
Set-ConsolidatedText $p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange "Slide 3"

# Slide 5 title: "Slide" " " "4" -> "Slide 4"
Set-ConsolidatedText $p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange "Slide 4"

# Slide 6 title: "Slide" " " "5" -> "Slide 5"
Set-ConsolidatedText $p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange "Slide 5"

# Slide 7 (the blank slide) speaker notes: many single-word runs ->
# one consolidated run.
$notes7 = $p.Slides.Item(7).NotesPage
Set-ConsolidatedText $notes7.Shapes.Item(2).TextFrame.TextRange "This is a blank slide: does it have a footer?"
